# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a handful of
# leve rows across the crafting job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2000520
$ws.Range("I38").Value = 200
$ws.Range("J38").Value = 5001000
$ws.Range("K38").Value = 600
$ws.Range("L38").Value = 15003000
$ws.Range("M38").Value = -228
$ws.Range("N38").Value = -15003744

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1278.5714
$ws.Range("I2").Value = 1262.5
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 1262.5
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = -1149.5
$ws.Range("N2").Value = -1526

$ws.Range("H32").Value = 9054.857
$ws.Range("I32").Value = 10084.7705
$ws.Range("J32").Value = 5759.1333
$ws.Range("K32").Value = 10084.7705
$ws.Range("L32").Value = 5759.1333
$ws.Range("M32").Value = -9797.770500000001
$ws.Range("N32").Value = -6333.1333

$ws.Range("H45").Value = 2717.8667
$ws.Range("I45").Value = 2704.9092
$ws.Range("J45").Value = 2753.5
$ws.Range("K45").Value = 2704.9092
$ws.Range("L45").Value = 2753.5
$ws.Range("M45").Value = -2327.9092
$ws.Range("N45").Value = -3507.5

$ws.Range("H61").Value = 9805989
$ws.Range("I61").Value = 12501960
$ws.Range("J61").Value = 2455.818
$ws.Range("K61").Value = 12501960
$ws.Range("L61").Value = 2455.818
$ws.Range("M61").Value = -12501748
$ws.Range("N61").Value = -2879.818

$ws.Range("H74").Value = 7144120.5
$ws.Range("I74").Value = 8475470
$ws.Range("J74").Value = 3246.182
$ws.Range("K74").Value = 8475470
$ws.Range("L74").Value = 3246.182
$ws.Range("M74").Value = -8474596
$ws.Range("N74").Value = -4994.182

$ws.Range("H77").Value = 7144120.5
$ws.Range("I77").Value = 8475470
$ws.Range("J77").Value = 3246.182
$ws.Range("K77").Value = 42377350
$ws.Range("L77").Value = 16230.91
$ws.Range("M77").Value = -42372982
$ws.Range("N77").Value = -24966.91

$ws.Range("H102").Value = 10250
$ws.Range("I102").Value = 9571.429
$ws.Range("K102").Value = 9571.429
$ws.Range("M102").Value = -7949.429

$ws.Range("H116").Value = 1278.5714
$ws.Range("I116").Value = 1262.5
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 1262.5
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = 1031.5
$ws.Range("N116").Value = -5888

$ws.Range("H132").Value = 5816018
$ws.Range("I132").Value = 10418291
$ws.Range("J132").Value = 2620.5264
$ws.Range("K132").Value = 31254873
$ws.Range("L132").Value = 7861.5792
$ws.Range("M132").Value = -31252343
$ws.Range("N132").Value = -12921.5792

$ws.Range("H135").Value = 61610
$ws.Range("J135").Value = 61610
$ws.Range("L135").Value = 61610
$ws.Range("N135").Value = -71750

$ws.Range("H136").Value = 9805989
$ws.Range("I136").Value = 12501960
$ws.Range("J136").Value = 2455.818
$ws.Range("K136").Value = 37505880
$ws.Range("L136").Value = 7367.454000000001
$ws.Range("M136").Value = -37503330
$ws.Range("N136").Value = -12467.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1278.5714
$ws.Range("I3").Value = 1262.5
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 1262.5
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = -1148.5
$ws.Range("N3").Value = -1528

$ws.Range("H99").Value = 2077.5
$ws.Range("I99").Value = 1924
$ws.Range("J99").Value = 2333.3333
$ws.Range("K99").Value = 1924
$ws.Range("L99").Value = 2333.3333
$ws.Range("M99").Value = -426
$ws.Range("N99").Value = -5329.3333

$ws.Range("H105").Value = 3378.2083
$ws.Range("I105").Value = 2382.8572
$ws.Range("J105").Value = 3788.0588
$ws.Range("K105").Value = 2382.8572
$ws.Range("L105").Value = 3788.0588
$ws.Range("M105").Value = -635.8571999999999
$ws.Range("N105").Value = -7282.0588

$ws.Range("H134").Value = 5154.873
$ws.Range("I134").Value = 4317.2324
$ws.Range("J134").Value = 6955.8
$ws.Range("K134").Value = 12951.6972
$ws.Range("L134").Value = 20867.4
$ws.Range("M134").Value = -10416.6972
$ws.Range("N134").Value = -25937.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 338.64102
$ws.Range("I107").Value = 314.54544
$ws.Range("J107").Value = 369.82352
$ws.Range("K107").Value = 314.54544
$ws.Range("L107").Value = 369.82352
$ws.Range("M107").Value = 1605.45456
$ws.Range("N107").Value = -4209.82352

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1395.5483
$ws.Range("J113").Value = 1619.7646
$ws.Range("L113").Value = 4859.293799999999
$ws.Range("N113").Value = -9199.293799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1811.6383
$ws.Range("I102").Value = 2658.875
$ws.Range("J102").Value = 927.56525
$ws.Range("K102").Value = 2658.875
$ws.Range("L102").Value = 927.56525
$ws.Range("M102").Value = -1036.875
$ws.Range("N102").Value = -4171.56525

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 4870.7334
$ws.Range("I132").Value = 4041.4167
$ws.Range("J132").Value = 5423.6113
$ws.Range("K132").Value = 12124.2501
$ws.Range("L132").Value = 16270.8339
$ws.Range("M132").Value = -9594.250100000001
$ws.Range("N132").Value = -21330.8339

$ws.Range("H141").Value = 373964.5
$ws.Range("J141").Value = 373964.5
$ws.Range("L141").Value = 373964.5
$ws.Range("N141").Value = -384324.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 943.3333
$ws.Range("I46").Value = 843.3333
$ws.Range("J46").Value = 1043.3334
$ws.Range("K46").Value = 843.3333
$ws.Range("L46").Value = 1043.3334
$ws.Range("M46").Value = -655.3333
$ws.Range("N46").Value = -1419.3334

$ws.Range("H61").Value = 1868.8
$ws.Range("I61").Value = 1758.8
$ws.Range("J61").Value = 1978.8
$ws.Range("K61").Value = 1758.8
$ws.Range("L61").Value = 1978.8
$ws.Range("M61").Value = -1556.8
$ws.Range("N61").Value = -2382.8

$ws.Range("H113").Value = 1868.8
$ws.Range("I113").Value = 1758.8
$ws.Range("J113").Value = 1978.8
$ws.Range("K113").Value = 1758.8
$ws.Range("L113").Value = 1978.8
$ws.Range("M113").Value = 411.2
$ws.Range("N113").Value = -6318.8

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = -71750

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
